$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("D4").Value = "2016-01-21 02:20:04"
$wsZhCn.Range("G4").Value = "2016-01-21 02:20:54"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("D4").Value = "2016-01-21 02:20:16"
$wsDeDe.Range("G4").Value = "2016-01-21 02:21:14"
